$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 212 (pushes old rows 212..244 down to 213..245)
$ws.Rows.Item(212).Insert()

# Populate the new row 212 with the new data point
$ws.Cells.Item(212, 1).Value = 7
$ws.Cells.Item(212, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(212, 3).Value = "Ñuble"
$ws.Cells.Item(212, 4).Value = 44776
$ws.Cells.Item(212, 5).Value = 16
$ws.Cells.Item(212, 6).Value = 100112043
$ws.Cells.Item(212, 7).Value = "Pepino ensalada"
$ws.Cells.Item(212, 8).Value = "Sin especificar"
$ws.Cells.Item(212, 9).Value = "Primera"
$ws.Cells.Item(212, 10).Value = 100
$ws.Cells.Item(212, 11).Value = 19000
$ws.Cells.Item(212, 12).Value = 20000
$ws.Cells.Item(212, 13).Value = 19500
$ws.Cells.Item(212, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(212, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(212, 16).Value = 325
$ws.Cells.Item(212, 17).Value = 60
$ws.Cells.Item(212, 18).Value = "Hortaliza"
